# Apply the "Add 2020-12-09 data" update to the Fonds de solidarite
# (volet 2, regional x categorie juridique) sheet.
#
# For each touched row we update:
#   - column C (nombre_aides)  -> new count
#   - column D (montant_total) -> new amount
#
# The source data stores every value as text (inline strings), even the
# numeric-looking ones, so each new value is written with a leading
# apostrophe to force Excel to keep storing it as text rather than
# silently re-typing the cell as a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;   C = "1437";  D = "9595127.25" },
    @{ Row = 6;   C = "1020";  D = "5886590.53" },
    @{ Row = 24;  C = "202";   D = "1090826.23" },
    @{ Row = 51;  C = "1143";  D = "8291040.22" },
    @{ Row = 52;  C = "794";   D = "5030000.28" },
    @{ Row = 55;  C = "10144"; D = "29325153.25" },
    @{ Row = 61;  C = "6684";  D = "28881781.83" },
    @{ Row = 96;  C = "655";   D = "4458316.04" },
    @{ Row = 101; C = "1497";  D = "3795586.09" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = "'" + $u.C
    $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
}
